# Commit: "Add files via upload"
#
# The shared string "A,4PC=1PC" is replaced everywhere by "AB,4PC=1PC".
# In the workbook this string only ever appears in the "Unit_Rule_Agg"
# column (E) for rows 44-57 of 工作表1, so updating that range reproduces
# the same data change as the upstream diff (the many other <v> index
# shifts in the diff are just Excel's automatic shared-string
# renumbering after the old string is dropped and the new one appended -
# they are not independent edits).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E44:E57").Value = "AB,4PC=1PC"

# Match the author's final selection: the bottom (frozen-below) pane now
# has E45:E57 selected with E45 as the active cell.
[void]$ws.Range("E45:E57").Select()
